$wb = $excel.ActiveWorkbook

# "html" sheet (3rd sheet) gets a new first column and new rows of content.
$ws = $wb.Worksheets.Item(3)

# Insert a new column at A, shifting the existing A:H content to B:I.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column.
$ws.Range("A1").Value = "Controller"

# New "Patient" sub-section header at row 16, styled with the red fill
# (same fill used elsewhere in the sheet, but without the centered
# vertical alignment that the other red-filled cells use).
$ws.Range("A16").Value = "Patient"
$ws.Range("A16").Interior.Color = 255

# Row 17: new
$ws.Range("B17").Value = "new"
$ws.Range("B17").Interior.Color = 255
$ws.Range("B17").VerticalAlignment = -4108
$ws.Range("F17").Value = "yes"
$ws.Range("H17").Value = "no"
$ws.Range("I17").Value = "first argument in fomr cannot contain or nill or be empty this might help http://stackoverflow.com/questions/17635634/first-argument-in-form-cannot-contain-nil-or-be-empty-comments"

# Row 18: edit
$ws.Range("B18").Value = "edit"
$ws.Range("B18").Interior.Color = 255
$ws.Range("B18").VerticalAlignment = -4108
$ws.Range("F18").Value = "yes"
$ws.Range("H18").Value = "no"
$ws.Range("I18").Value = "undefined method name when rendering this page"
$ws.Range("J18").Value = "Content of the page not tested, only rending"

# Row 19: show
$ws.Range("B19").Value = "show"
$ws.Range("B19").Interior.Color = 255
$ws.Range("B19").VerticalAlignment = -4108
$ws.Range("F19").Value = "yes"
$ws.Range("H19").Value = "no"
$ws.Range("I19").Value = "undefined method name when rendering this page"
$ws.Range("J19").Value = "Content of the page not tested, only rending"

# Row 20: _form
$ws.Range("B20").Value = "_form"
$ws.Range("B20").Interior.Color = 255
$ws.Range("B20").VerticalAlignment = -4108
$ws.Range("F20").Value = "yes"
$ws.Range("H20").Value = "no"
$ws.Range("I20").Value = "first argument in fomr cannot contain or nill or be empty this might help http://stackoverflow.com/questions/17635634/first-argument-in-form-cannot-contain-nil-or-be-empty-comments"

# Make "html" the active/selected sheet and select B20 (matches the
# saved selection state in the workbook).
$ws.Activate()
[void]$ws.Range("B20").Select()
